# Auto-generated edit script: applies the "gh-pages output @456a3b4" update
# to sheet1 ("展览") and sheet4 ("全部类型").
$wb = $excel.ActiveWorkbook

function Set-StyledIndexCell($ws, $row, $col, $value) {
    # Recreate the bold/centered/thin-border "index column" style (style 1)
    # that every A-column data cell in these sheets uses, without leaving
    # the cell pointing at a freshly-minted duplicate style.
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

function Set-PlainTextCell($ws, $row, $col, $value) {
    # Force literal text storage (so date-shaped strings like "2024-08-17"
    # are not reinterpreted as date serials), then drop back to the default
    # "Normal" style so no stray number-format style lingers on the cell.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$ws = $wb.Worksheets.Item(1)

# --- plain vote-count (column F) bumps, no row movement ---
$ws.Cells.Item(4, 6).Value = 7966
$ws.Cells.Item(5, 6).Value = 98
$ws.Cells.Item(9, 6).Value = 114
$ws.Cells.Item(10, 6).Value = 465
$ws.Cells.Item(13, 6).Value = 451
$ws.Cells.Item(14, 6).Value = 68
$ws.Cells.Item(15, 6).Value = 75
$ws.Cells.Item(17, 6).Value = 5853
$ws.Cells.Item(18, 6).Value = 182
$ws.Cells.Item(19, 6).Value = 263
$ws.Cells.Item(20, 6).Value = 1790

# --- insert 2 new rows at 21; this pushes the old
#     "銀魂" row (21) down to 23, and the old
#     "SSS" row (22) down to 24 ---
$ws.Range("A21:A22").EntireRow.Insert()

# Row 21: in-place replacement -- old "銀魂" event becomes the
#     new "...水千丞周边预约票" event; A21 (index 20) is unchanged.
Set-StyledIndexCell $ws 21 1 20
Set-PlainTextCell $ws 21 2 "2024-08-17"
$ws.Cells.Item(21, 3).Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws.Cells.Item(21, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Cells.Item(21, 5).Value = "2024.08.17 09:30-08.17 17:00"
$ws.Cells.Item(21, 6).Value = 2
$ws.Cells.Item(21, 7).Value = 0.1
$ws.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

# Row 22: brand-new "...水千丞签售预约票" event.
Set-StyledIndexCell $ws 22 1 21
Set-PlainTextCell $ws 22 2 "2024-08-17"
$ws.Cells.Item(22, 3).Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws.Cells.Item(22, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Cells.Item(22, 5).Value = "2024.08.17 09:30-08.17 17:00"
$ws.Cells.Item(22, 6).Value = 5
$ws.Cells.Item(22, 7).Value = 0.1
$ws.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws.Cells.Item(22, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

# Row 23: the original "銀魂" row, shifted down intact by the insert;
#     only its sequence index (A23) and vote count (F23) change.
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 6).Value = 238

# Row 24: the original "SSS" row, shifted down intact by the insert;
#     only its sequence index (A24) and vote count (F24) change.
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 6).Value = 390

# ---- Sheet 4: "全部类型" (all types) ----
$ws = $wb.Worksheets.Item(4)

# --- plain vote-count (column F) bumps, no row movement ---
$ws.Cells.Item(4, 6).Value = 7966
$ws.Cells.Item(5, 6).Value = 98
$ws.Cells.Item(9, 6).Value = 114
$ws.Cells.Item(10, 6).Value = 465
$ws.Cells.Item(13, 6).Value = 451
$ws.Cells.Item(14, 6).Value = 68
$ws.Cells.Item(15, 6).Value = 75
$ws.Cells.Item(18, 6).Value = 5853
$ws.Cells.Item(20, 6).Value = 182
$ws.Cells.Item(21, 6).Value = 263
$ws.Cells.Item(22, 6).Value = 1790

# --- insert 2 new rows at 23; this pushes the old
#     "銀魂" row (23) down to 25, and the old
#     "SSS" row (24) down to 26 ---
$ws.Range("A23:A24").EntireRow.Insert()

# Row 23: in-place replacement -- old "銀魂" event becomes the
#     new "...水千丞周边预约票" event; A23 (index 22) is unchanged.
Set-StyledIndexCell $ws 23 1 22
Set-PlainTextCell $ws 23 2 "2024-08-17"
$ws.Cells.Item(23, 3).Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws.Cells.Item(23, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Cells.Item(23, 5).Value = "2024.08.17 09:30-08.17 17:00"
$ws.Cells.Item(23, 6).Value = 2
$ws.Cells.Item(23, 7).Value = 0.1
$ws.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws.Cells.Item(23, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

# Row 24: brand-new "...水千丞签售预约票" event.
Set-StyledIndexCell $ws 24 1 23
Set-PlainTextCell $ws 24 2 "2024-08-17"
$ws.Cells.Item(24, 3).Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws.Cells.Item(24, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Cells.Item(24, 5).Value = "2024.08.17 09:30-08.17 17:00"
$ws.Cells.Item(24, 6).Value = 5
$ws.Cells.Item(24, 7).Value = 0.1
$ws.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws.Cells.Item(24, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

# Row 25: the original "銀魂" row, shifted down intact by the insert;
#     only its sequence index (A25) and vote count (F25) change.
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 6).Value = 238

# Row 26: the original "SSS" row, shifted down intact by the insert;
#     only its sequence index (A26) and vote count (F26) change.
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 6).Value = 390

Write-Output "edit complete"
